$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 4 swap for columns D, L, M, N, O, P, S

$ws.Range("D2").Value = 44881
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11250
$ws.Range("O2").Value = 11250
$ws.Range("P2").Value = 11250
$ws.Range("S2").Value = 11250

$ws.Range("D4").Value = 44923
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7625
$ws.Range("S4").Value = 7625
